$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.529.28"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.833.15"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9975"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.65"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9936"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4468"
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.54"
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07790"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.144"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.34"
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9937"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.362"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.549"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "1.832.68"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.49"
$ws.Range("E17").Value = "  +13.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001089"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06369"
$ws.Range("E19").Value = "  -5.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9952"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.67"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.395"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5398"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").Value = "28.576.07"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.88"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.244"
$ws.Range("E26").Value = "  -8.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.99"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.05"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.397"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").Value = "2.040.95"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "130.21"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.228"
$ws.Range("E32").Value = "  -4.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.897"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09294"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.672"
$ws.Range("E35").Value = "  -7.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.90"
$ws.Range("E36").Value = "  +5.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02370"
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2207"
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6675"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.222"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06284"
$ws.Range("E41").Value = "  -1.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.197"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.134"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9933"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.401"
$ws.Range("E45").Value = "  -3.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.95"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6133"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.772"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.83"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.049"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.02"
$ws.Range("E51").Value = "  +1.84%  "
